$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1596385542168675
$ws.Range("C2").Value = 0.6144578313253012
$ws.Range("J2").Value = 0.01807228915662651
$ws.Range("P2").Value = 0.1204819277108434
$ws.Range("S2").Value = 0.08734939759036145
$ws.Range("B3").Value = 0.004694835680751174
$ws.Range("C3").Value = 0.04225352112676056
$ws.Range("J3").Value = 0.02816901408450704
$ws.Range("P3").Value = 0.7370892018779343
$ws.Range("S3").Value = 0.1877934272300469
$ws.Range("J4").Value = 0.01886792452830189
$ws.Range("P4").Value = 0.6415094339622641
$ws.Range("S4").Value = 0.3396226415094339
$ws.Range("B6").Value = 0.05857740585774059
$ws.Range("D6").Value = 0.02928870292887029
$ws.Range("F6").Value = 0.09205020920502092
$ws.Range("J6").Value = 0.2343096234309623
$ws.Range("O6").Value = 0.03765690376569038
$ws.Range("Q6").Value = 0.1548117154811715
$ws.Range("R6").Value = 0.05439330543933055
$ws.Range("S6").Value = 0.3389121338912134
$ws.Range("B7").Value = 0.1204188481675393
$ws.Range("D7").Value = 0.02617801047120419
$ws.Range("E7").Value = 0.005235602094240838
$ws.Range("F7").Value = 0.03141361256544502
$ws.Range("J7").Value = 0.1151832460732984
$ws.Range("O7").Value = 0.02617801047120419
$ws.Range("Q7").Value = 0.1518324607329843
$ws.Range("R7").Value = 0.07853403141361257
$ws.Range("S7").Value = 0.4450261780104712
$ws.Range("B8").Value = 0.1058495821727019
$ws.Range("D8").Value = 0.008356545961002786
$ws.Range("F8").Value = 0.07520891364902507
$ws.Range("J8").Value = 0.1002785515320334
$ws.Range("O8").Value = 0.01114206128133705
$ws.Range("Q8").Value = 0.1587743732590529
$ws.Range("R8").Value = 0.1086350974930362
$ws.Range("S8").Value = 0.4317548746518106
$ws.Range("B9").Value = 0.1153846153846154
$ws.Range("D9").Value = 0.01282051282051282
$ws.Range("F9").Value = 0.06837606837606838
$ws.Range("J9").Value = 0.1367521367521368
$ws.Range("O9").Value = 0.0170940170940171
$ws.Range("Q9").Value = 0.1752136752136752
$ws.Range("R9").Value = 0.08974358974358974
$ws.Range("S9").Value = 0.3846153846153846
$ws.Range("B10").Value = 0.132629992464205
$ws.Range("D10").Value = 0.02863602110022608
$ws.Range("E10").Value = 0.0007535795026375283
$ws.Range("F10").Value = 0.0746043707611153
$ws.Range("J10").Value = 0.1243406179351922
$ws.Range("O10").Value = 0.01582516955538809
$ws.Range("Q10").Value = 0.2155237377543331
$ws.Range("R10").Value = 0.08138658628485305
$ws.Range("S10").Value = 0.3262999246420497
$ws.Range("G11").Value = 0.1346801346801347
$ws.Range("J11").Value = 0.07407407407407407
$ws.Range("K11").Value = 0.1750841750841751
$ws.Range("L11").Value = 0.6026936026936027
$ws.Range("S11").Value = 0.01346801346801347
$ws.Range("G12").Value = 0.6983240223463687
$ws.Range("J12").Value = 0.2290502793296089
$ws.Range("K12").Value = 0.0111731843575419
$ws.Range("L12").Value = 0.0111731843575419
$ws.Range("S12").Value = 0.05027932960893855
$ws.Range("G13").Value = 0.7317073170731707
$ws.Range("J13").Value = 0.2439024390243902
$ws.Range("S13").Value = 0.02439024390243903
$ws.Range("F15").Value = 0.0371900826446281
$ws.Range("H15").Value = 0.1239669421487603
$ws.Range("I15").Value = 0.06611570247933884
$ws.Range("J15").Value = 0.359504132231405
$ws.Range("K15").Value = 0.05785123966942149
$ws.Range("M15").Value = 0.02066115702479339
$ws.Range("O15").Value = 0.08264462809917356
$ws.Range("S15").Value = 0.2520661157024793
$ws.Range("F16").Value = 0.0179372197309417
$ws.Range("H16").Value = 0.1121076233183857
$ws.Range("I16").Value = 0.1255605381165919
$ws.Range("J16").Value = 0.3901345291479821
$ws.Range("K16").Value = 0.1345291479820628
$ws.Range("M16").Value = 0.02242152466367713
$ws.Range("O16").Value = 0.05829596412556054
$ws.Range("S16").Value = 0.1390134529147982
$ws.Range("F17").Value = 0.008928571428571428
$ws.Range("H17").Value = 0.1473214285714286
$ws.Range("I17").Value = 0.1138392857142857
$ws.Range("J17").Value = 0.4709821428571428
$ws.Range("K17").Value = 0.08928571428571429
$ws.Range("M17").Value = 0.01116071428571429
$ws.Range("O17").Value = 0.06026785714285714
$ws.Range("S17").Value = 0.09821428571428571
$ws.Range("F18").Value = 0.01538461538461539
$ws.Range("H18").Value = 0.1435897435897436
$ws.Range("I18").Value = 0.1128205128205128
$ws.Range("J18").Value = 0.4205128205128205
$ws.Range("K18").Value = 0.09230769230769231
$ws.Range("M18").Value = 0.01538461538461539
$ws.Range("N18").Value = 0.005128205128205128
$ws.Range("O18").Value = 0.08205128205128205
$ws.Range("S18").Value = 0.1128205128205128
$ws.Range("F19").Value = 0.01848874598070739
$ws.Range("H19").Value = 0.1728295819935691
$ws.Range("I19").Value = 0.09163987138263666
$ws.Range("J19").Value = 0.3826366559485531
$ws.Range("K19").Value = 0.1117363344051447
$ws.Range("M19").Value = 0.01688102893890675
$ws.Range("N19").Value = 0.0008038585209003215
$ws.Range("O19").Value = 0.07877813504823152
$ws.Range("S19").Value = 0.1262057877813505
